$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets("ALC")
$ws.Range("H32").Value = 2774.0625
$ws.Range("I32").Value = 1903
$ws.Range("J32").Value = 3645.125
$ws.Range("K32").Value = 1903
$ws.Range("L32").Value = 3645.125
$ws.Range("M32").Value = -1577
$ws.Range("N32").Value = -4297.125
$ws.Range("H33").Value = 330.82608
$ws.Range("I33").Value = 339.0476
$ws.Range("K33").Value = 339.0476
$ws.Range("M33").Value = -110.0476
$ws.Range("H58").Value = 1237.6666
$ws.Range("I58").Value = 463
$ws.Range("K58").Value = 1389
$ws.Range("M58").Value = -1239
$ws.Range("H74").Value = 9500.714
$ws.Range("J74").Value = 9233.333000000001
$ws.Range("L74").Value = 9233.333000000001
$ws.Range("N74").Value = -11105.333
$ws.Range("H77").Value = 9500.714
$ws.Range("J77").Value = 9233.333000000001
$ws.Range("L77").Value = 46166.665
$ws.Range("N77").Value = -55526.665
$ws.Range("H86").Value = 1350.8572
$ws.Range("J86").Value = 1610.6
$ws.Range("L86").Value = 1610.6
$ws.Range("N86").Value = -3856.6
$ws.Range("H89").Value = 1350.8572
$ws.Range("J89").Value = 1610.6
$ws.Range("L89").Value = 8053
$ws.Range("N89").Value = -19285
$ws.Range("H92").Value = 50549.85
$ws.Range("I92").Value = 53131.42
$ws.Range("K92").Value = 53131.42
$ws.Range("M92").Value = -51883.42
$ws.Range("H100").Value = 2633.25
$ws.Range("I100").Value = 2198.4614
$ws.Range("J100").Value = 4517.3335
$ws.Range("K100").Value = 2198.4614
$ws.Range("L100").Value = 4517.3335
$ws.Range("M100").Value = -1657.4614
$ws.Range("N100").Value = -5599.3335
$ws.Range("H129").Value = 1841.75
$ws.Range("I129").Value = 819.7143
$ws.Range("J129").Value = 3272.6
$ws.Range("K129").Value = 2459.1429
$ws.Range("L129").Value = 9817.799999999999
$ws.Range("M129").Value = 2540.8571
$ws.Range("N129").Value = -19817.8
$ws.Range("H137").Value = 3118.6191
$ws.Range("I137").Value = 1250.5
$ws.Range("J137").Value = 3558.1765
$ws.Range("K137").Value = 3751.5
$ws.Range("L137").Value = 10674.5295
$ws.Range("M137").Value = -1201.5
$ws.Range("N137").Value = -15774.5295
$ws.Range("H138").Value = 1721.25
$ws.Range("I138").Value = 1605
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 4815
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = 325
$ws.Range("N138").Value = -19280
$ws.Range("H141").Value = 25999.75
$ws.Range("I141").Value = 25999.75
$ws.Range("K141").Value = 77999.25
$ws.Range("M141").Value = -72819.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets("ARM")
$ws.Range("H2").Value = 469.5
$ws.Range("I2").Value = 469.5
$ws.Range("K2").Value = 469.5
$ws.Range("M2").Value = -356.5
$ws.Range("H32").Value = 1589.3334
$ws.Range("I32").Value = 1441.0435
$ws.Range("K32").Value = 1441.0435
$ws.Range("M32").Value = -1154.0435
$ws.Range("H61").Value = 6955.85
$ws.Range("I61").Value = 6169.769
$ws.Range("K61").Value = 6169.769
$ws.Range("M61").Value = -5957.769
$ws.Range("H74").Value = 1948
$ws.Range("I74").Value = 1948
$ws.Range("K74").Value = 1948
$ws.Range("M74").Value = -1074
$ws.Range("H77").Value = 1948
$ws.Range("I77").Value = 1948
$ws.Range("K77").Value = 9740
$ws.Range("M77").Value = -5372
$ws.Range("H97").Value = 1545.4117
$ws.Range("I97").Value = 1102.2858
$ws.Range("J97").Value = 3613.3333
$ws.Range("K97").Value = 1102.2858
$ws.Range("L97").Value = 3613.3333
$ws.Range("M97").Value = -606.2858000000001
$ws.Range("N97").Value = -4605.3333
$ws.Range("H110").Value = 9320
$ws.Range("I110").Value = 11014.4
$ws.Range("J110").Value = 848
$ws.Range("K110").Value = 11014.4
$ws.Range("L110").Value = 848
$ws.Range("M110").Value = -8969.4
$ws.Range("N110").Value = -4938
$ws.Range("H116").Value = 469.5
$ws.Range("I116").Value = 469.5
$ws.Range("K116").Value = 469.5
$ws.Range("M116").Value = 1824.5
$ws.Range("H122").Value = 1600
$ws.Range("I122").Value = 1600
$ws.Range("K122").Value = 4800
$ws.Range("M122").Value = -2350
$ws.Range("H132").Value = 3890.9167
$ws.Range("I132").Value = 3177
$ws.Range("K132").Value = 9531
$ws.Range("M132").Value = -7001
$ws.Range("H136").Value = 6955.85
$ws.Range("I136").Value = 6169.769
$ws.Range("K136").Value = 18509.307
$ws.Range("M136").Value = -15959.307

# ---- Sheet: BSM ----
$ws = $wb.Worksheets("BSM")
$ws.Range("H3").Value = 469.5
$ws.Range("I3").Value = 469.5
$ws.Range("K3").Value = 469.5
$ws.Range("M3").Value = -355.5
$ws.Range("H80").Value = 5593.357
$ws.Range("J80").Value = 9505.875
$ws.Range("L80").Value = 9505.875
$ws.Range("N80").Value = -11501.875
$ws.Range("H83").Value = 5593.357
$ws.Range("J83").Value = 9505.875
$ws.Range("L83").Value = 47529.375
$ws.Range("N83").Value = -57513.375
$ws.Range("H86").Value = 5854.0967
$ws.Range("I86").Value = 2277.5789
$ws.Range("J86").Value = 11516.917
$ws.Range("K86").Value = 2277.5789
$ws.Range("L86").Value = 11516.917
$ws.Range("M86").Value = -1154.5789
$ws.Range("N86").Value = -13762.917
$ws.Range("H89").Value = 5854.0967
$ws.Range("I89").Value = 2277.5789
$ws.Range("J89").Value = 11516.917
$ws.Range("K89").Value = 11387.8945
$ws.Range("L89").Value = 57584.585
$ws.Range("M89").Value = -5771.8945
$ws.Range("N89").Value = -68816.58499999999
$ws.Range("H105").Value = 3600.8
$ws.Range("I105").Value = 3929.7144
$ws.Range("J105").Value = 2833.3333
$ws.Range("K105").Value = 3929.7144
$ws.Range("L105").Value = 2833.3333
$ws.Range("M105").Value = -2182.7144
$ws.Range("N105").Value = -6327.3333
$ws.Range("H134").Value = 5089.6665
$ws.Range("I134").Value = 4816.7856
$ws.Range("K134").Value = 14450.3568
$ws.Range("M134").Value = -11915.3568

# ---- Sheet: CRP ----
$ws = $wb.Worksheets("CRP")
$ws.Range("H22").Value = 3636930.5
$ws.Range("I22").Value = 505
$ws.Range("J22").Value = 10000675
$ws.Range("K22").Value = 505
$ws.Range("L22").Value = 10000675
$ws.Range("M22").Value = -155
$ws.Range("N22").Value = -10001375
$ws.Range("H31").Value = 2077.875
$ws.Range("I31").Value = 1838.9375
$ws.Range("J31").Value = 2555.75
$ws.Range("K31").Value = 1838.9375
$ws.Range("L31").Value = 2555.75
$ws.Range("M31").Value = -1543.9375
$ws.Range("N31").Value = -3145.75
$ws.Range("H34").Value = 2077.875
$ws.Range("I34").Value = 1838.9375
$ws.Range("J34").Value = 2555.75
$ws.Range("K34").Value = 1838.9375
$ws.Range("L34").Value = 2555.75
$ws.Range("M34").Value = -1636.9375
$ws.Range("N34").Value = -2959.75
$ws.Range("H58").Value = 4839.75
$ws.Range("I58").Value = 2634.875
$ws.Range("K58").Value = 2634.875
$ws.Range("M58").Value = -2431.875
$ws.Range("H122").Value = 6022.1113
$ws.Range("I122").Value = 5964.5884
$ws.Range("K122").Value = 17893.7652
$ws.Range("M122").Value = -15443.7652
$ws.Range("H136").Value = 4839.75
$ws.Range("I136").Value = 2634.875
$ws.Range("K136").Value = 7904.625
$ws.Range("M136").Value = -5354.625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets("CUL")
$ws.Range("H36").Value = 91
$ws.Range("I36").Value = 91
$ws.Range("K36").Value = 273
$ws.Range("M36").Value = -104
$ws.Range("H132").Value = 2165.5
$ws.Range("I132").Value = 1831.3334
$ws.Range("J132").Value = 2499.6667
$ws.Range("K132").Value = 16482.0006
$ws.Range("L132").Value = 22497.0003
$ws.Range("M132").Value = -13952.0006
$ws.Range("N132").Value = -27557.0003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets("GSM")
$ws.Range("H122").Value = 2782.2778
$ws.Range("I122").Value = 2944.4614
$ws.Range("K122").Value = 8833.3842
$ws.Range("M122").Value = -6383.3842
$ws.Range("H132").Value = 2322.25
$ws.Range("I132").Value = 2322.25
$ws.Range("K132").Value = 6966.75
$ws.Range("M132").Value = -4436.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets("LTW")
$ws.Range("H22").Value = 3206.3333
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3206.3333
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3206.3333
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -3796.3333
$ws.Range("H27").Value = 3206.3333
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 3206.3333
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 3206.3333
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -3420.3333
$ws.Range("H100").Value = 2407.3333
$ws.Range("I100").Value = 2315.6667
$ws.Range("K100").Value = 2315.6667
$ws.Range("M100").Value = -1774.6667
$ws.Range("H122").Value = 3284.75
$ws.Range("I122").Value = 3284.75
$ws.Range("K122").Value = 9854.25
$ws.Range("M122").Value = -7404.25
$ws.Range("H132").Value = 2142.5862
$ws.Range("I132").Value = 2052.3
$ws.Range("K132").Value = 6156.900000000001
$ws.Range("M132").Value = -3626.900000000001
$ws.Range("H136").Value = 3369.4348
$ws.Range("I136").Value = 1849.9
$ws.Range("K136").Value = 5549.700000000001
$ws.Range("M136").Value = -2999.700000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets("WVR")
$ws.Range("H45").Value = 47969.69
$ws.Range("I45").Value = 37434.832
$ws.Range("K45").Value = 37434.832
$ws.Range("M45").Value = -36943.832
$ws.Range("H122").Value = 4089.4707
$ws.Range("I122").Value = 4108.7144
$ws.Range("K122").Value = 12326.1432
$ws.Range("M122").Value = -9876.143199999999
$ws.Range("H136").Value = 10109.6875
$ws.Range("I136").Value = 11670.926
$ws.Range("J136").Value = 1679
$ws.Range("K136").Value = 35012.778
$ws.Range("L136").Value = 5037
$ws.Range("M136").Value = -32462.778
$ws.Range("N136").Value = -10137
